# Weekly update: insert the latest Fruta/Mango price record for
# "Vega Monumental Concepción" as a new row right before the existing
# row 58, pushing the rest of the table (old rows 58-159) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 58 - this shifts rows 58:159 down to 59:160
# and Excel auto-extends the used range / dimension accordingly.
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new weekly record.
$ws.Cells.Item(58, 1).Value  = 11
$ws.Cells.Item(58, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(58, 3).Value  = "Bíobío"
$ws.Cells.Item(58, 4).Value  = 45036
$ws.Cells.Item(58, 5).Value  = 8
$ws.Cells.Item(58, 6).Value  = "Fruta"
$ws.Cells.Item(58, 7).Value  = 100108
$ws.Cells.Item(58, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(58, 9).Value  = 100108002
$ws.Cells.Item(58, 10).Value = "Mango"
$ws.Cells.Item(58, 11).Value = "Sin especificar"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 200
$ws.Cells.Item(58, 14).Value = 7000
$ws.Cells.Item(58, 15).Value = 7500
$ws.Cells.Item(58, 16).Value = 7200
$ws.Cells.Item(58, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(58, 18).Value = "Perú"
$ws.Cells.Item(58, 19).Value = 1800
$ws.Cells.Item(58, 20).Value = 4
